$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (pushes current rows 11..127 down to 12..128)
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with this week's data (matches the
# constant columns used throughout the sheet for this market/product).
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Terminal La Palmera de La Serena"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44630
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 100112044
$ws.Range("G11").Value = "Perejil"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 2500
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 2750
$ws.Range("N11").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 1833
$ws.Range("Q11").Value = 1.5
$ws.Range("R11").Value = "Hortaliza"
